$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new record (PMC3749898 / Grudnikoff et al.) ---
$ws.Range("B4").Value = 'JOUR'
$ws.Range("C4").Value = 'Objective: Behavioral and psychiatric disorders are common in youth with rapid-onset obesity with hypothalamic dysfunction, hypoventilation, and autonomic dysregulation (ROHHAD). We outline a rational approach to psychiatric treatment of a patient with a complex medical condition. Method(s): We report the course of symptoms in a teen with ROHHAD, the inpatient treatment, and review current evidence for use of psychopharmacologic agents in youth with sleep and anxiety disturbances. Result(s): A 14-year-old female began rapidly gaining weight as a preschooler, developed hormonal imbalance, and mixed sleep apnea. Consultation was requested after a month of ROHHAD exacerbation, with severe anxiety, insomnia, and auditory hallucinations. Olanzapine and citalopram were helpful in controlling the symptoms. Following discharge, the patient gained weight and olanzapine was discontinued. Lorazepam was started in coordination with pulmonary service. Relevant pharmacologic considerations included risk of respiratory suppression, history of paradoxical reaction to hypnotics, hepatic isoenzyme interactions and side effects of antipsychotics. Conclusion(s): Core symptoms of ROHHAD may precipitate psychiatric disorders. A systematic evidence-based approach to psychopharmacology is necessary in the setting of psychiatric consultation.'
$ws.Range("D4").Value = 'E. Grudnikoff, Department of Psychiatry, Zucker Hillside Hospital, North Shore-LIJ Health System, New York, NY, United States. E-mail: egrudnikof@nshs.edu'
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '3'
$ws.Range("F4").ClearFormats()
$ws.Range("G4").Value = 'adolescent, aggression, *anxiety disorder/dt [Drug Therapy], *anxiety disorder/si [Side Effect], anxiety disorder/dt [Drug Therapy], anxiety disorder/si [Side Effect], article, auditory hallucination/dt [Drug Therapy], *autonomic dysfunction, avoidance behavior, bedtime dosage, behavior change, case report, developmental disorder/dt [Drug Therapy], differential diagnosis, disease association, disease exacerbation, disease severity, distress syndrome/si [Side Effect], drug dose increase, drug safety, drug tolerability, drug withdrawal, electrolyte disturbance/dt [Drug Therapy], endocrine disease, evening dosage, female, growth disorder/dt [Drug Therapy], human, insomnia/dt [Drug Therapy], irritability, medical history, nightmare/dt [Drug Therapy], nightmare/si [Side Effect], onset age, risk assessment, side effect/si [Side Effect], sleep disorder/dt [Drug Therapy], sleep disorder/si [Side Effect], sleep disordered breathing, vivid dream/si [Side Effect], weight gain, citalopram/dt [Drug Therapy], desmopressin/dt [Drug Therapy], estrogen/dt [Drug Therapy], eszopiclone/ae [Adverse Drug Reaction], eszopiclone/cb [Drug Combination], eszopiclone/dt [Drug Therapy], growth hormone/dt [Drug Therapy], lorazepam, melatonin/ae [Adverse Drug Reaction], melatonin/cb [Drug Combination], melatonin/dt [Drug Therapy], olanzapine/ae [Adverse Drug Reaction], olanzapine/dt [Drug Therapy], zolpidem/ae [Adverse Drug Reaction], zolpidem/cb [Drug Combination], zolpidem/dt [Drug Therapy], *nocturnal anxiety, *rapid onset obesity hypothalamic dysfunction hypoventilation and autonomic dysregulation'
$ws.Range("H4").Value = 'English'
$ws.Range("J4").Value = '1719-8429'
$ws.Range("K4").Value = 'Nocturnal anxiety in a youth with rapid-onset obesity, hypothalamic dysfunction, hypoventilation, and autonomic dysregulation (rohhad)'
$ws.Range("L4").Value = 'Journal of the Canadian Academy of Child and Adolescent Psychiatry'
$ws.Range("M4").Value = 'Nocturnal anxiety in a youth with rapid-onset obesity, hypothalamic dysfunction, hypoventilation, and autonomic dysregulation (rohhad)'
$ws.Range("N4").Value = 'http://www.cacap-acpea.org/uploads/documents//Nocturnal_Anxiety_Grudnikoff.pdfhttp://ovidsp.ovid.com/ovidweb.cgi?T=JS&PAGE=reference&D=emed14&NEWS=N&AN=369579073'
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = '22'
$ws.Range("O4").ClearFormats()
$ws.Range("P4").Value = 23450
$ws.Range("S4").Value = 'Grudnikoff, Eugene Foley, Carmel Poole, Claudette Theodosiadis, Eva eng Canada J Can Acad Child Adolesc Psychiatry. 2013 Aug;22(3):235-7.'
$ws.Range("T4").Value = 'Grudnikoff, Eugene, Foley, Carmel, Theodosiadis, Eva, Poole, Claudette'
$ws.Range("Y4").Value = 23970913
$ws.Range("Z4").Value = 'PMC3749898'
$ws.Range("AA4").Value = 'Aug'
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = '2013-08-24'
$ws.Range("AB4").ClearFormats()
$ws.Range("AC4").Value = '235-237'
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 1618

# --- A4: record_id, bold + bordered + centered (matches header/id-column styling) ---
$ws.Range("A4").Value = 18527
$ws.Range("A4").Font.Name = "Calibri"
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4160

Write-Host "Row 4 populated."
